# "updated UI for batch operation"
# - Rename the "Not Started" activation-status value to the machine-friendly
#   "not_started" token used by the batch-operation backend.
# - Leave the sheet with cell E13 selected/active (reflects the UI state at
#   the time of the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update all "Not Started" activation status cells (column C, data rows) to
# the new lowercase/underscored status token "not_started".
$ws.Range("C2:C496").Replace("Not Started", "not_started")

# Reflect the active cell/selection shown in the sheet view.
$ws.Cells.Item(13, 5).Select()
